$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.613.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.595.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.820.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.584.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.595.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.277.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.772'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.732.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.103'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.91%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
